$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $text
}

# Data row starting at table row 1 (original: 15÷9=,11÷2=,88÷7=,30÷4=,35÷4=)
Set-CellText $t 1 1 "12÷8="
Set-CellText $t 1 2 "56÷5="
Set-CellText $t 1 3 "41÷9="
Set-CellText $t 1 4 "62÷3="
Set-CellText $t 1 5 "81÷9="

# Data row starting at table row 5 (original: 68÷4=,91÷9=,10÷7=,87÷9=,90÷4=)
Set-CellText $t 5 1 "43÷9="
Set-CellText $t 5 2 "44÷8="
Set-CellText $t 5 3 "38÷6="
Set-CellText $t 5 4 "25÷8="
Set-CellText $t 5 5 "21÷7="

# Data row starting at table row 9 (original: 52÷4=,59÷5=,15÷9=,65÷5=,79÷4=)
Set-CellText $t 9 1 "35÷3="
Set-CellText $t 9 2 "20÷5="
Set-CellText $t 9 3 "84÷7="
Set-CellText $t 9 4 "33÷2="
Set-CellText $t 9 5 "90÷4="

# Data row starting at table row 13 (original: 11÷3=,42÷2=,13÷2=,31÷6=,12÷9=)
Set-CellText $t 13 1 "22÷7="
Set-CellText $t 13 2 "29÷6="
Set-CellText $t 13 3 "14÷7="
Set-CellText $t 13 4 "43÷7="
Set-CellText $t 13 5 "63÷7="

# Data row starting at table row 17 (original: 42÷5=,25÷9=,20÷6=,44÷4=,62÷6=)
Set-CellText $t 17 1 "65÷2="
Set-CellText $t 17 2 "33÷5="
Set-CellText $t 17 3 "13÷6="
Set-CellText $t 17 4 "50÷8="
Set-CellText $t 17 5 "52÷3="

Write-Host "Done updating cells."
